$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 through 9 (Flight_Smoke, Flight_Regression, Hotel_Smoke, Hotel_Regression, Regression_All)
$ws.Range("A5:B9").EntireRow.Delete()

# Rename remaining module values to strip the "_All" suffix
$ws.Range("A2").Value = "Flight"
$ws.Range("A3").Value = "Hotel"
$ws.Range("A4").Value = "Holiday"

# Update the selection to match the new active cell
$ws.Range("B4").Select()
